# Add 9 new "geohash_1".."geohash_9" field rows to the "Fields" sheet,
# right after the "geoform" row (old row 27) and before the "id" row
# (old row 28), per the commit "Added geohash fields to SearchMetadata
# documentation."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fields")

# Insert 9 blank rows starting at row 28 (pushes old rows 28-92 down to 37-101).
$ws.Rows("28:36").Insert()

$desc = "An encoded string that represents the geographic coordinates of the centroid of a spatial extent. This can be used for searching and plotting."
$src  = "Encoded from centroid of northBoundCoord, sourthBoundCoord, eastBoundCoord, westBoundCoord"

# Columns A-E first (matches original authoring order of the shared-string table)...
for ($i = 1; $i -le 9; $i++) {
  $r = 27 + $i
  $ws.Cells.Item($r, 1).Value2 = "geohash_$i"
  $ws.Cells.Item($r, 2).Value2 = "string"
  $ws.Cells.Item($r, 3).Value2 = "Yes"
  $ws.Cells.Item($r, 4).Value2 = "S"
  $ws.Cells.Item($r, 5).Value2 = $desc
  $ws.Rows.Item($r).RowHeight = 45
}

# ...then column F (the "Encoded from ..." note) in a second pass.
for ($i = 1; $i -le 9; $i++) {
  $r = 27 + $i
  $ws.Cells.Item($r, 6).Value2 = $src
}

# The conditional formatting over column A covered A2:A93; after inserting
# 9 rows it should cover A2:A102.
$fcs = $ws.Range("A2").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
  $fcs.Item($i).ModifyAppliesToRange($ws.Range("A2:A102"))
}

# The hidden _FilterDatabase defined name tracked the sheet's filter range;
# keep it in sync with the new used range (A1:H101).
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
  $n = $names.Item($i)
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=Fields!`$A`$1:`$H`$101"
  }
}

# The EML sheet picked up an explicit portrait page setup in this revision.
$wsEml = $wb.Worksheets.Item("EML")
$wsEml.PageSetup.Orientation = 1

# Re-select the cell the author ended up on and bring "Fields" to front
# (it was tab 4 "EML" before the edit).
$ws.Activate()
$ws.Range("E30").Select()
